$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 907, pushing the existing row 907 (and everything
# below it) down by one. This matches the diff: a "2026/03/01" / 日 row
# (time 4, ranking 201) is added right after the existing 2026/03/01 / 日
# row, and the rest of the table (2026/12/29 ... 2027/01/05) shifts from
# rows 907-948 to rows 908-949.
$ws.Rows(907).Insert()

# Column A holds dates stored as plain text (e.g. "2026/03/01"), not real
# date serials. Force text formatting before assigning so Excel doesn't
# auto-convert the string into a date value, then restore the default
# "Normal" style so the new cell doesn't pick up a stray number format
# that the rest of the column doesn't have.
$ws.Range("A907").NumberFormat = "@"
$ws.Range("A907").Value = "2026/03/01"
$ws.Range("A907").Style = "Normal"

$ws.Range("B907").Value = "日"
$ws.Range("C907").Value = 4
$ws.Range("D907").Value = 201
